$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "74.711.31"
$ws.Range("E2").Value = "  +8.72%  "

# Row 3
$ws.Range("D3").Value = "2.602.97"
$ws.Range("E3").Value = "  +7.25%  "

# Row 4
$ws.Range("D4").Value = "'1.00"
$ws.Range("E4").Value = "  +0.16%  "

# Row 5
$ws.Range("D5").Value = "'185.75"
$ws.Range("E5").Value = "  +15.63%  "

# Row 6
$ws.Range("D6").Value = "'584.02"
$ws.Range("E6").Value = "  +4.53%  "

# Row 7
$ws.Range("E7").Value = "  +0.00%  "

# Row 8
$ws.Range("B8").Value = "XRP"
$ws.Range("C8").Value = "https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp"
$ws.Range("D8").Value = "'0.539"
$ws.Range("E8").Value = "  +5.80%  "

# Row 9
$ws.Range("B9").Value = "Dogecoin"
$ws.Range("C9").Value = "https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge"
$ws.Range("D9").Value = "'0.208"
$ws.Range("E9").Value = "  +23.79%  "

# Row 10
$ws.Range("D10").Value = "2.602.76"
$ws.Range("E10").Value = "  +7.32%  "

# Row 11
$ws.Range("E11").Value = "  +0.03%  "

# Row 12
$ws.Range("D12").Value = "'0.364"
$ws.Range("E12").Value = "  +10.08%  "

# Row 13
$ws.Range("E13").Value = "  +4.98%  "

# Row 14
$ws.Range("D14").Value = "'0.0000192"
$ws.Range("E14").Value = "  +9.58%  "

# Row 15
$ws.Range("D15").Value = "74.661.55"
$ws.Range("E15").Value = "  +8.87%  "

# Row 16
$ws.Range("D16").Value = "3.087.89"

# Row 17
$ws.Range("D17").Value = "'26.29"
$ws.Range("E17").Value = "  +13.61%  "

# Row 18
$ws.Range("D18").Value = "2.613.50"
$ws.Range("E18").Value = "  +7.67%  "

# Row 19
$ws.Range("D19").Value = "'9.31"
$ws.Range("E19").Value = "  +34.84%  "

# Row 20
$ws.Range("D20").Value = "'11.80"
$ws.Range("E20").Value = "  +12.62%  "

# Row 21
$ws.Range("D21").Value = "'377.50"
$ws.Range("E21").Value = "  +12.75%  "

# Row 22
$ws.Range("D22").Value = "'2.29"
$ws.Range("E22").Value = "  +18.37%  "

# Row 23
$ws.Range("D23").Value = "'4.09"
$ws.Range("E23").Value = "  +7.40%  "

# Row 24
$ws.Range("D24").Value = "'1.00"
$ws.Range("E24").Value = "  +0.14%  "

# Row 25
$ws.Range("D25").Value = "'70.31"
$ws.Range("E25").Value = "  +5.03%  "

# Row 26
$ws.Range("D26").Value = "'4.20"
$ws.Range("E26").Value = "  +14.30%  "

# Row 27
$ws.Range("D27").Value = "'9.32"
$ws.Range("E27").Value = "  +14.09%  "

# Row 28
$ws.Range("D28").Value = "2.741.37"
$ws.Range("E28").Value = "  +7.32%  "

# Row 29
$ws.Range("D29").Value = "'0.993"
$ws.Range("E29").Value = "  -0.70%  "

# Row 30
$ws.Range("D30").Value = "0.0₃0953"
$ws.Range("E30").Value = "  +16.69%  "

# Row 31
$ws.Range("D31").Value = "'1.39"
$ws.Range("E31").Value = "  +20.61%  "

# Row 32
$ws.Range("D32").Value = "'7.96"
$ws.Range("E32").Value = "  +11.86%  "

# Row 33
$ws.Range("D33").Value = "'508.80"
$ws.Range("E33").Value = "  +19.15%  "

# Row 34
$ws.Range("D34").Value = "'1.75"
$ws.Range("E34").Value = "  +8.80%  "

# Row 35
$ws.Range("D35").Value = "'1.00"
$ws.Range("E35").Value = "  +0.17%  "

# Row 36
$ws.Range("D36").Value = "'0.122"
$ws.Range("E36").Value = "  +15.78%  "

# Row 37
$ws.Range("D37").Value = "'159.08"
$ws.Range("E37").Value = "  -0.95%  "

# Row 38
$ws.Range("D38").Value = "'19.27"
$ws.Range("E38").Value = "  +7.69%  "

# Row 39
$ws.Range("D39").Value = "'19.37"
$ws.Range("E39").Value = "  +1.84%  "

# Row 41
$ws.Range("B41").Value = "Stacks"
$ws.Range("C41").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D41").Value = "'1.70"
$ws.Range("E41").Value = "  +13.42%  "

# Row 42
$ws.Range("B42").Value = "RenderToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/vfo5XYwcV+rendertoken-render"
$ws.Range("D42").Value = "'4.91"
$ws.Range("E42").Value = "  +13.32%  "

# Row 43
$ws.Range("E43").Value = "  +10.29%  "

# Row 44
$ws.Range("D44").Value = "'2.44"
$ws.Range("E44").Value = "  +19.63%  "

# Row 45
$ws.Range("D45").Value = "'157.22"
$ws.Range("E45").Value = "  +19.27%  "

# Row 46
$ws.Range("E46").Value = "  +10.26%  "

# Row 47
$ws.Range("E47").Value = "  +4.01%  "

# Row 48
$ws.Range("D48").Value = "'0.0842"
$ws.Range("E48").Value = "  +17.43%  "

# Row 49
$ws.Range("D49").Value = "'3.63"
$ws.Range("E49").Value = "  +8.82%  "

# Row 50
$ws.Range("D50").Value = "'0.526"
$ws.Range("E50").Value = "  +9.39%  "

# Row 51
$ws.Range("D51").Value = "'20.04"
$ws.Range("E51").Value = "  +18.89%  "
